$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(1,1).Value = "TEST123"
$v = $ws.Cells.Item(1,1).Value2
Write-Output "val: $v"
